# Fix bug where the springer API would not properly fill in the abstract
# in some situations. Refresh Abstract (column D) and Authors (column E)
# text for the affected rows with the corrected, fully-populated values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "[Ali%Dadgari%NULL%1,            Tengku%Aizan Hamid%NULL%1,            Mohammad Nazrul%Hakim%NULL%1,            Reza%Chaman%NULL%1,            Seyed Abbas%Mousavi%NULL%1,            Lim%Poh Hin%NULL%1,            Leila%Dadvar%NULL%1]"
$ws.Range("D3").Value = "`nAims`nThe purpose of this study is to investigate whether an in-home, individually tailored intervention is efficacious in promoting increases in physical activity (PA) and improvements in physical functioning (PF) in low-active older adults.`nMethods`nParticipants were randomized to two groups for the 8-week intervention.`n The enhanced physical activity (EPA) group received individualized exercise programming, including personalized step goals and a resistance band training program, and the standard of care (SoC) group received a general activity goal.`n Pre- and post-intervention PF measures included choice step reaction time, knee extension/flexion strength, hand grip strength, and 8 ft up and go test completion time.`nResults`nThirty-nine subjects completed this study (74.6 ± 6.4 years).`n Significant increases in steps/day were observed for both the EPA and SoC groups, although the improvements in the EPA group were significantly higher when including only those who adhered to weekly step goals.`n Both groups experienced significant PF improvements, albeit greater in the EPA group for the 8 ft up and go test and knee extension strength.`nConclusion`nA low cost, in-home intervention elicited improvements in both PA and PF.`n Future research is warranted to expand upon the size and scope of this study, exploring dose thresholds (and time frames) for PA to improve PF and strategies to further bolster adherence rates to maximize intervention benefits.`n"
$ws.Range("E3").Value = "[ Christopher J.%Dondzila%null%0,  Ann M.%Swartz%null%1,  Kevin G.%Keenan%null%1,  Amy E.%Harley%null%1,  Razia%Azen%null%1,  Scott J.%Strath%null%1]"
$ws.Range("E4").Value = "[Ryoichi%Ema%NULL%1,            Shunsuke%Ohki%NULL%1,            Hirokazu%Takayama%NULL%1,            Yuji%Kobayashi%NULL%1,            Ryota%Akagi%NULL%1]"
$ws.Range("E5").Value = "[Tsung-Jen%Hsieh%tjhsieh@nhri.edu.tw%1,            Shin-Chang%Su%maxmax0220@hotmail.com%1,            Chun-Wei%Chen%chunwei.peter@gmail.com%1,            Yaw-Wen%Kang%castle.kang@gmail.com%1,            Ming-Hsia%Hu%mhh@ntu.edu.tw%1,            Li-Lin%Hsu%llh@firdi.org.tw%1,            Szu-Yun%Wu%s.wu@ibms.sinica.edu.tw%1,            Likwang%Chen%likwang@nhri.edu.tw%1,            Hsing-Yi%Chang%hsingyi@nhri.edu.tw%1,            Shao-Yuan%Chuang%chuangsy@nhri.edu.tw%1,            Wen-Harn%Pan%pan@ibms.sinica.edu.tw%1,            Chih-Cheng%Hsu%cch@nhri.edu.tw%1]"
$ws.Range("E6").Value = "[Steve%Iliffe%xref no email%1,    Denise%Kendrick%xref no email%1,    Richard%Morris%xref no email%1,    Mark%Griffin%xref no email%1,    Deborah%Haworth%xref no email%1,    Hannah%Carpenter%xref no email%1,    Tahir%Masud%xref no email%1,    Dawn A%Skelton%xref no email%1,    Susie%Dinan-Young%xref no email%1,    Ann%Bowling%xref no email%1,    Heather%Gage%xref no email%1]"
$ws.Range("E7").Value = "[Nicole%Kahle%xref no email%1,    Michael A.%Tevald%xref no email%1]"
$ws.Range("E9").Value = "[Andr\u00e9%Lacroix%xref no email%1,    Reto W.%Kressig%xref no email%1,    Thomas%Muehlbauer%xref no email%1,    Yves J.%Gschwind%xref no email%1,    Barbara%Pfenninger%xref no email%1,    Othmar%Bruegger%xref no email%1,    Urs%Granacher%xref no email%1]"
$ws.Range("E10").Value = "[Ahamed%Yasmin%coreGivesNoEmail%1,          Close%Jacqueline%coreGivesNoEmail%1,          Cook%Wendy L.%coreGivesNoEmail%1,          Donaldson%Meghan G.%coreGivesNoEmail%1,          Graf%Peter%coreGivesNoEmail%1,          Khan%Karim M.%coreGivesNoEmail%1,          Liu-Ambrose%Teresa%coreGivesNoEmail%1,          Lord%Stephen R.%coreGivesNoEmail%1]"
$ws.Range("E11").Value = "[Kohei%Maruya%NULL%1,            Yasuyoshi%Asakawa%NULL%1,            Hideaki%Ishibashi%NULL%1,            Hiroaki%Fujita%NULL%1,            Tomoyuki%Arai%NULL%1,            Haruyasu%Yamaguchi%NULL%1]"
$ws.Range("E12").Value = "[M. E.%Nelson%xref no email%1,    J. E.%Layne%xref no email%1,    M. J.%Bernstein%xref no email%1,    A.%Nuernberger%xref no email%1,    C.%Castaneda%xref no email%1,    D.%Kaliton%xref no email%1,    J.%Hausdorff%xref no email%1,    J. O.%Judge%xref no email%1,    D. M.%Buchner%xref no email%1,    R.%Roubenoff%xref no email%1,    M. A.%Fiatarone Singh%xref no email%1]"
$ws.Range("E14").Value = "[Oliver J.%Perkin%o.j.perkin@bath.ac.uk%1,            Polly M.%McGuigan%NULL%2,            Polly M.%McGuigan%NULL%0,            Keith A.%Stokes%NULL%1]"
$ws.Range("E15").Value = "[Maria%Tsekoura%NULL%1,            Evdokia%Billis%NULL%1,            Elias%Tsepis%NULL%1,            Zacharias%Dimitriadis%NULL%1,            Charalampos%Matzaroglou%NULL%1,            Minos%Tyllianakis%NULL%1,            Elias%Panagiotopoulos%NULL%1,            John%Gliatis%NULL%1]"
$ws.Range("D16").Value = "`n                Background and aims: Home-based exercise is a viable solution for frail elderly individuals with difficulties in reaching exercise facilities outside home.`n The aim of this study was to determine the effects of a home-based video exercise program on physiological performance, functional capacity and health-related quality of life.`n Methods: Community-dwelling frail women ≥75 yrs, receiving public home care, were randomized into a training group (n=30) and a control group (n=31).`n Participants exercised for 26 minutes, three times per week for five months.`n Both groups received a bi-weekly telephone call.`n The effect of intervention was evaluated by the physical performance test, mobility-tiredness score, maximal isometric handgrip and biceps strength, lower limb explosive power, repeated chair rise (5 times), 10-m maximal walking-speed, semi-tandem balance, and health-related quality of life, as measured by EQ-5D and self-rated health.`n Results: Twenty-five participants (83%) in the training group and 28 (90%) in the control group completed the project.`n Adherence to the training protocol was on average 89.2%.`n At follow-up, between-group analysis revealed a significant difference only in EQ-5D (valued by time-trade-off tariffs), resulting from a significant decrease observed in the control group and a trend towards an increase in the training group (p=0.082).`n Significant within-group improvements, ranging from 8–35%, were also observed for the physical performance test, mobility-tiredness score, handgrip, biceps strength, chair rise, and 10-m maximal walking-speed in the training group, and for walking-speed and self-rated health in the control group.`n Conclusions: These results suggest that home-based training for frail older women using an exercise video induces lasting health-related quality-of-life (EQ-5D).`n In addition, a tendency towards improvements in physiological performance and functional capacity was observed.`n"
$ws.Range("E16").Value = "[ Sonja%Vestergaard%null%0,  Christian%Kronborg%null%1,  Lis%Puggaard%null%1]"
